$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("H1").Value = "Tempo Heuristica"
$ws.Range("I1").Value = "Tempo Total"

# Row 2 updates
$ws.Range("C2").Value = 2380
$ws.Range("F2").Value = "23 -> 22 -> 21 -> 67 -> 68 -> 70 -> 72 -> 10 -> 7 -> 4 -> 1 -> 2 -> 5 -> 8 -> 12 -> 16 -> 17 -> 20 -> 24 -> 32 -> 44 -> 46 -> 45 -> 43 -> 32 -> 24 -> 23"
$ws.Range("G2").Value = 0.02661943435668945
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.02661943435668945

# Row 3 updates
$ws.Range("C3").Value = 3091
$ws.Range("F3").Value = "48 -> 47 -> 46 -> 45 -> 49 -> 53 -> 55 -> 57 -> 60 -> 61 -> 64 -> 65 -> 66 -> 22 -> 21 -> 18 -> 14 -> 15 -> 14 -> 10 -> 7 -> 4 -> 1 -> 77 -> 78 -> 79 -> 81 -> 3 -> 25 -> 26 -> 27 -> 28 -> 29 -> 30 -> 37 -> 38 -> 41 -> 44 -> 46 -> 47 -> 48"
$ws.Range("G3").Value = 0.01713728904724121
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.01713728904724121

# Row 4 updates
$ws.Range("C4").Value = 2750
$ws.Range("E4").Value = "11 -> 34 -> 31 -> 40 -> 63 -> 64 -> 11"
$ws.Range("F4").Value = "11 -> 12 -> 13 -> 29 -> 28 -> 34 -> 28 -> 29 -> 30 -> 31 -> 38 -> 39 -> 40 -> 42 -> 41 -> 32 -> 24 -> 23 -> 22 -> 66 -> 65 -> 64 -> 61 -> 62 -> 63 -> 62 -> 61 -> 64 -> 65 -> 66 -> 22 -> 21 -> 18 -> 14 -> 10 -> 11"
$ws.Range("G4").Value = 0.01559281349182129
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.01559281349182129
